# "Server to Local on 09 02 24" - refresh the Assessment Report sheet with the
# latest export: drop the unused "Section2" column and replace the single
# sample row with the current set of participant results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assessment Report")

# ---------------------------------------------------------------------------
# 1. Remove the "Section2" column (old column D). This shifts Secured Marks,
#    Percentage, Result and Remarks one column to the left, matching the new
#    A1:G12 layout.
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).Delete()

# ---------------------------------------------------------------------------
# 2. Replace the existing data row (row 2) and append the rest of the rows
#    from the latest server export.
#    Columns after the delete: A Submitted Date, B Participant Name,
#    C Section1, D Secured Marks, E Percentage, F Result, G Remarks.
# ---------------------------------------------------------------------------
$submitted = Get-Date -Year 2024 -Month 2 -Day 8 -Hour 11 -Minute 52 -Second 42

# name, Section1, SecuredMarks, Percentage(text), Result, Remarks
$data = @(
    @("Bhakata Ram Suna",        25,  25,  "25",  "Not Cleared", "Auto Submitted Assessment (Moved Outside the assessment boundary)."),
    @("Swati Swarupa Rajguru",   30,  30,  "30",  "Not Cleared", "Auto Submitted Assessment (Moved Outside the assessment boundary)."),
    @("Sanigdha Mohanty",        29,  29,  "29",  "Not Cleared", "Auto Submitted Assessment (Moved Outside the assessment boundary)."),
    @("Satyajeet behera",        63,  63,  "63",  "Not Cleared", "Auto Submitted Assessment (Moved Outside the assessment boundary)."),
    @("Satyajeet behera",         0,   0,  "0",   "Not Cleared", "Switching of Tab's detected."),
    @("Smitaranjan Samantaray",  84,  84,  "84",  "Not Cleared", "N/A"),
    @("Soumyaranjan Sethy",     100, 100,  "100", "Cleared",     "N/A"),
    @("Jagadish Prasad Dash",    61,  61,  "61",  "Not Cleared", "Auto Submitted Assessment (Moved Outside the assessment boundary)."),
    @("Bhakata Ram Suna",        76,  76,  "76",  "Not Cleared", "N/A"),
    @("Gagan Kumar Behura",      57,  57,  "57",  "Not Cleared", "Time Out Auto Submitting the Assessment."),
    @("Kollu Nagarjuna",         65,  65,  "65",  "Not Cleared", "N/A")
)

# Carry the date style (m/d/yyyy, same as the existing A2) down through A12
# without introducing a new number-format style entry.
$ws.Range("A2").Copy()
$ws.Range("A3:A12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $submitted
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]

    # Column E ("Percentage") holds a numeric-looking value stored as TEXT in
    # the source export. Build it with TEXT() in a scratch cell and paste
    # only the resulting value across so the destination keeps a plain/default
    # cell style (matches the un-styled <c t="s"> cells in the export).
    $helper = $ws.Cells.Item($r, 50)
    $helper.Formula = '=TEXT(' + $row[1] + ',"0")'
    $helper.Copy()
    $ws.Cells.Item($r, 5).PasteSpecial(-4163) | Out-Null
    $helper.Delete(-4159)
}
$excel.CutCopyMode = 0
